$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 40000
$ws.Range("J3").Value = 40000
$ws.Range("L3").Value = 40000
$ws.Range("N3").Value = -40228
$ws.Range("H9").Value = 764230.9
$ws.Range("J9").Value = 591.8570999999999
$ws.Range("L9").Value = 591.8570999999999
$ws.Range("N9").Value = -929.8570999999999
$ws.Range("H12").Value = 150
$ws.Range("I12").Value = 150
$ws.Range("K12").Value = 150
$ws.Range("M12").Value = 20
$ws.Range("H19").Value = 2891.5557
$ws.Range("I19").Value = 2366.3333
$ws.Range("J19").Value = 3416.7778
$ws.Range("K19").Value = 2366.3333
$ws.Range("L19").Value = 3416.7778
$ws.Range("M19").Value = -2191.3333
$ws.Range("N19").Value = -3766.7778
$ws.Range("H43").Value = 19439.3
$ws.Range("I43").Value = 22842
$ws.Range("J43").Value = 11499.667
$ws.Range("K43").Value = 22842
$ws.Range("L43").Value = 11499.667
$ws.Range("M43").Value = -22773
$ws.Range("N43").Value = -11637.667
$ws.Range("H64").Value = 7849.25
$ws.Range("I64").Value = 7899.5
$ws.Range("J64").Value = 7799
$ws.Range("K64").Value = 7899.5
$ws.Range("L64").Value = 7799
$ws.Range("M64").Value = -7651.5
$ws.Range("N64").Value = -8295
$ws.Range("H67").Value = 7849.25
$ws.Range("I67").Value = 7899.5
$ws.Range("J67").Value = 7799
$ws.Range("K67").Value = 7899.5
$ws.Range("L67").Value = 7799
$ws.Range("M67").Value = -7041.5
$ws.Range("N67").Value = -9515
$ws.Range("H102").Value = 40000
$ws.Range("J102").Value = 40000
$ws.Range("L102").Value = 40000
$ws.Range("N102").Value = -46490
$ws.Range("H112").Value = 3814.111
$ws.Range("I112").Value = 1399.5
$ws.Range("K112").Value = 4198.5
$ws.Range("M112").Value = -3090.5
$ws.Range("H113").Value = 4994
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 2043.5
$ws.Range("I4").Value = 387.5
$ws.Range("J4").Value = 3699.5
$ws.Range("K4").Value = 387.5
$ws.Range("L4").Value = 3699.5
$ws.Range("M4").Value = -271.5
$ws.Range("N4").Value = -3931.5
$ws.Range("H32").Value = 1071.3164
$ws.Range("I32").Value = 756.0540999999999
$ws.Range("K32").Value = 756.0540999999999
$ws.Range("M32").Value = -469.0540999999999
$ws.Range("H45").Value = 8425.799999999999
$ws.Range("I45").Value = 12638.315
$ws.Range("J45").Value = 1149.6364
$ws.Range("K45").Value = 12638.315
$ws.Range("L45").Value = 1149.6364
$ws.Range("M45").Value = -12261.315
$ws.Range("N45").Value = -1903.6364
$ws.Range("H61").Value = 8278.823
$ws.Range("I61").Value = 8362.833000000001
$ws.Range("K61").Value = 8362.833000000001
$ws.Range("M61").Value = -8150.833000000001
$ws.Range("H110").Value = 3701.1667
$ws.Range("I110").Value = 2735
$ws.Range("J110").Value = 4667.3335
$ws.Range("K110").Value = 2735
$ws.Range("L110").Value = 4667.3335
$ws.Range("M110").Value = -690
$ws.Range("N110").Value = -8757.333500000001
$ws.Range("H136").Value = 8278.823
$ws.Range("I136").Value = 8362.833000000001
$ws.Range("K136").Value = 25088.499
$ws.Range("M136").Value = -22538.499

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2997.5
$ws.Range("J80").Value = 2997.5
$ws.Range("L80").Value = 2997.5
$ws.Range("N80").Value = -4993.5
$ws.Range("H83").Value = 2997.5
$ws.Range("J83").Value = 2997.5
$ws.Range("L83").Value = 14987.5
$ws.Range("N83").Value = -24971.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 2000
$ws.Range("I23").Value = 2000
$ws.Range("K23").Value = 2000
$ws.Range("M23").Value = -1760
$ws.Range("H27").Value = 2000
$ws.Range("I27").Value = 2000
$ws.Range("K27").Value = 2000
$ws.Range("M27").Value = -1808
$ws.Range("H31").Value = 5338.95
$ws.Range("I31").Value = 5891.143
$ws.Range("J31").Value = 4050.5
$ws.Range("K31").Value = 5891.143
$ws.Range("L31").Value = 4050.5
$ws.Range("M31").Value = -5596.143
$ws.Range("N31").Value = -4640.5
$ws.Range("H34").Value = 5338.95
$ws.Range("I34").Value = 5891.143
$ws.Range("J34").Value = 4050.5
$ws.Range("K34").Value = 5891.143
$ws.Range("L34").Value = 4050.5
$ws.Range("M34").Value = -5689.143
$ws.Range("N34").Value = -4454.5
$ws.Range("H58").Value = 4913.0527
$ws.Range("I58").Value = 5216.3228
$ws.Range("J58").Value = 3570
$ws.Range("K58").Value = 5216.3228
$ws.Range("L58").Value = 3570
$ws.Range("M58").Value = -5013.3228
$ws.Range("N58").Value = -3976
$ws.Range("H99").Value = 6747
$ws.Range("I99").Value = 6205.1665
$ws.Range("K99").Value = 6205.1665
$ws.Range("M99").Value = -4707.1665
$ws.Range("H126").Value = 6747
$ws.Range("I126").Value = 6205.1665
$ws.Range("K126").Value = 18615.4995
$ws.Range("M126").Value = -16145.4995
$ws.Range("H132").Value = 4664.8
$ws.Range("I132").Value = 4480.8623
$ws.Range("J132").Value = 9999
$ws.Range("K132").Value = 13442.5869
$ws.Range("L132").Value = 29997
$ws.Range("M132").Value = -10912.5869
$ws.Range("N132").Value = -35057
$ws.Range("H134").Value = 7036.0415
$ws.Range("I134").Value = 6528.95
$ws.Range("K134").Value = 19586.85
$ws.Range("M134").Value = -17051.85
$ws.Range("H135").Value = 89900
$ws.Range("J135").Value = 89900
$ws.Range("L135").Value = 89900
$ws.Range("N135").Value = -100040
$ws.Range("H136").Value = 4913.0527
$ws.Range("I136").Value = 5216.3228
$ws.Range("J136").Value = 3570
$ws.Range("K136").Value = 15648.9684
$ws.Range("L136").Value = 10710
$ws.Range("M136").Value = -13098.9684
$ws.Range("N136").Value = -15810
$ws.Range("H141").Value = 68548
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 68548
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 68548
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -78908

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 1722
$ws.Range("I60").Value = 374.5
$ws.Range("J60").Value = 2800
$ws.Range("K60").Value = 1123.5
$ws.Range("L60").Value = 8400
$ws.Range("M60").Value = -872.5
$ws.Range("N60").Value = -8902
$ws.Range("H100").Value = 11248.75
$ws.Range("J100").Value = 11248.75
$ws.Range("L100").Value = 33746.25
$ws.Range("N100").Value = -35368.25
$ws.Range("H117").Value = 1625
$ws.Range("I117").Value = 645.8333
$ws.Range("J117").Value = 2212.5
$ws.Range("K117").Value = 1937.4999
$ws.Range("L117").Value = 6637.5
$ws.Range("M117").Value = 1504.5001
$ws.Range("N117").Value = -13521.5
$ws.Range("H138").Value = 11360
$ws.Range("I138").Value = 5666.6665
$ws.Range("J138").Value = 19900
$ws.Range("K138").Value = 16999.9995
$ws.Range("L138").Value = 59700
$ws.Range("M138").Value = -11859.9995
$ws.Range("N138").Value = -69980

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2994.5
$ws.Range("I80").Value = 2990
$ws.Range("K80").Value = 2990
$ws.Range("M80").Value = -1992
$ws.Range("H83").Value = 2994.5
$ws.Range("I83").Value = 2990
$ws.Range("K83").Value = 14950
$ws.Range("M83").Value = -9958
$ws.Range("H102").Value = 2357.4546
$ws.Range("I102").Value = 2293.2
$ws.Range("K102").Value = 2293.2
$ws.Range("M102").Value = -671.1999999999998
$ws.Range("H132").Value = 21973
$ws.Range("I132").Value = 5188.1304
$ws.Range("J132").Value = 214999
$ws.Range("K132").Value = 15564.3912
$ws.Range("L132").Value = 644997
$ws.Range("M132").Value = -13034.3912
$ws.Range("N132").Value = -650057
$ws.Range("H133").Value = 83684.42999999999
$ws.Range("J133").Value = 83684.42999999999
$ws.Range("L133").Value = 83684.42999999999
$ws.Range("N133").Value = -93804.42999999999
$ws.Range("H138").Value = 79799.2
$ws.Range("J138").Value = 79799.2
$ws.Range("L138").Value = 79799.2
$ws.Range("N138").Value = -90079.2

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 9006.5
$ws.Range("I25").Value = 9006.5
$ws.Range("K25").Value = 9006.5
$ws.Range("M25").Value = -8776.5
$ws.Range("H55").Value = 272.4643
$ws.Range("J55").Value = 206.06667
$ws.Range("L55").Value = 206.06667
$ws.Range("N55").Value = -552.0666699999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 9797.5
$ws.Range("I52").Value = 9797.5
$ws.Range("K52").Value = 9797.5
$ws.Range("M52").Value = -9571.5
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H113").Value = 721.3570999999999
$ws.Range("J113").Value = 382.33334
$ws.Range("L113").Value = 1147.00002
$ws.Range("N113").Value = -5487.000019999999
